# Apply updated dSF ("F" column) values for several rows, per the
# "repull data, push all data, mean calculation" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value  = 0
$ws.Range("F5").Value  = 1
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = 7
$ws.Range("F15").Value = -1
$ws.Range("F16").Value = 0
